$d = $word.ActiveDocument

$d.Content.Find.Execute("40+21=61", $true, $false, $false, $false, $false, $true, 1, $false, "9-9=0", 2) | Out-Null
$d.Content.Find.Execute("24+70=94", $true, $false, $false, $false, $false, $true, 1, $false, "65-33=32", 2) | Out-Null
$d.Content.Find.Execute("48-20=28", $true, $false, $false, $false, $false, $true, 1, $false, "57+36=93", 2) | Out-Null
$d.Content.Find.Execute("64+7=71", $true, $false, $false, $false, $false, $true, 1, $false, "13+77=90", 2) | Out-Null
$d.Content.Find.Execute("12+7=19", $true, $false, $false, $false, $false, $true, 1, $false, "18+41=59", 2) | Out-Null
$d.Content.Find.Execute("64-25=39", $true, $false, $false, $false, $false, $true, 1, $false, "64-19=45", 2) | Out-Null
$d.Content.Find.Execute("98-38=60", $true, $false, $false, $false, $false, $true, 1, $false, "87-47=40", 2) | Out-Null
$d.Content.Find.Execute("45+6=51", $true, $false, $false, $false, $false, $true, 1, $false, "91-37=54", 2) | Out-Null
$d.Content.Find.Execute("87-9=78", $true, $false, $false, $false, $false, $true, 1, $false, "62-25=37", 2) | Out-Null
$d.Content.Find.Execute("36+49=85", $true, $false, $false, $false, $false, $true, 1, $false, "93-91=2", 2) | Out-Null
$d.Content.Find.Execute("1+31=32", $true, $false, $false, $false, $false, $true, 1, $false, "35-14=21", 2) | Out-Null
$d.Content.Find.Execute("92-77=15", $true, $false, $false, $false, $false, $true, 1, $false, "94-77=17", 2) | Out-Null
$d.Content.Find.Execute("73-33=40", $true, $false, $false, $false, $false, $true, 1, $false, "94-28=66", 2) | Out-Null
$d.Content.Find.Execute("44+34=78", $true, $false, $false, $false, $false, $true, 1, $false, "28-12=16", 2) | Out-Null
$d.Content.Find.Execute("67+26=93", $true, $false, $false, $false, $false, $true, 1, $false, "92-51=41", 2) | Out-Null
$d.Content.Find.Execute("18+42=60", $true, $false, $false, $false, $false, $true, 1, $false, "33+30=63", 2) | Out-Null
$d.Content.Find.Execute("87+0=87", $true, $false, $false, $false, $false, $true, 1, $false, "35+57=92", 2) | Out-Null
$d.Content.Find.Execute("15+11=26", $true, $false, $false, $false, $false, $true, 1, $false, "85-58=27", 2) | Out-Null
$d.Content.Find.Execute("13+2=15", $true, $false, $false, $false, $false, $true, 1, $false, "81-18=63", 2) | Out-Null
$d.Content.Find.Execute("22-0=22", $true, $false, $false, $false, $false, $true, 1, $false, "36-1=35", 2) | Out-Null
$d.Content.Find.Execute("41-28=13", $true, $false, $false, $false, $false, $true, 1, $false, "33-19=14", 2) | Out-Null
$d.Content.Find.Execute("45+27=72", $true, $false, $false, $false, $false, $true, 1, $false, "43+34=77", 2) | Out-Null
$d.Content.Find.Execute("85-55=30", $true, $false, $false, $false, $false, $true, 1, $false, "33+15=48", 2) | Out-Null
$d.Content.Find.Execute("23+13=36", $true, $false, $false, $false, $false, $true, 1, $false, "30+54=84", 2) | Out-Null
$d.Content.Find.Execute("64+14=78", $true, $false, $false, $false, $false, $true, 1, $false, "25+25=50", 2) | Out-Null
$d.Content.Find.Execute("16-16=0", $true, $false, $false, $false, $false, $true, 1, $false, "9+37=46", 2) | Out-Null
$d.Content.Find.Execute("35-18=17", $true, $false, $false, $false, $false, $true, 1, $false, "57+17=74", 2) | Out-Null
$d.Content.Find.Execute("67-50=17", $true, $false, $false, $false, $false, $true, 1, $false, "52-37=15", 2) | Out-Null
$d.Content.Find.Execute("27+68=95", $true, $false, $false, $false, $false, $true, 1, $false, "15-13=2", 2) | Out-Null
$d.Content.Find.Execute("82-40=42", $true, $false, $false, $false, $false, $true, 1, $false, "36+53=89", 2) | Out-Null
$d.Content.Find.Execute("92-86=6", $true, $false, $false, $false, $false, $true, 1, $false, "87-69=18", 2) | Out-Null
$d.Content.Find.Execute("72-45=27", $true, $false, $false, $false, $false, $true, 1, $false, "91-50=41", 2) | Out-Null
$d.Content.Find.Execute("32+25=57", $true, $false, $false, $false, $false, $true, 1, $false, "67+22=89", 2) | Out-Null
$d.Content.Find.Execute("99-50=49", $true, $false, $false, $false, $false, $true, 1, $false, "54+24=78", 2) | Out-Null
$d.Content.Find.Execute("46+26=72", $true, $false, $false, $false, $false, $true, 1, $false, "16+74=90", 2) | Out-Null
$d.Content.Find.Execute("48+17=65", $true, $false, $false, $false, $false, $true, 1, $false, "82-15=67", 2) | Out-Null
$d.Content.Find.Execute("22+40=62", $true, $false, $false, $false, $false, $true, 1, $false, "20+59=79", 2) | Out-Null
$d.Content.Find.Execute("84-42=42", $true, $false, $false, $false, $false, $true, 1, $false, "34+49=83", 2) | Out-Null
$d.Content.Find.Execute("55-9=46", $true, $false, $false, $false, $false, $true, 1, $false, "67+1=68", 2) | Out-Null
$d.Content.Find.Execute("41-0=41", $true, $false, $false, $false, $false, $true, 1, $false, "60+20=80", 2) | Out-Null
$d.Content.Find.Execute("44+26=70", $true, $false, $false, $false, $false, $true, 1, $false, "54+24=78", 2) | Out-Null
$d.Content.Find.Execute("44-10=34", $true, $false, $false, $false, $false, $true, 1, $false, "48-0=48", 2) | Out-Null
$d.Content.Find.Execute("78+3=81", $true, $false, $false, $false, $false, $true, 1, $false, "67-55=12", 2) | Out-Null
$d.Content.Find.Execute("62-41=21", $true, $false, $false, $false, $false, $true, 1, $false, "50-14=36", 2) | Out-Null
$d.Content.Find.Execute("80-9=71", $true, $false, $false, $false, $false, $true, 1, $false, "37+39=76", 2) | Out-Null
$d.Content.Find.Execute("10+71=81", $true, $false, $false, $false, $false, $true, 1, $false, "0+96=96", 2) | Out-Null
$d.Content.Find.Execute("49+9=58", $true, $false, $false, $false, $false, $true, 1, $false, "84-32=52", 2) | Out-Null
$d.Content.Find.Execute("0+52=52", $true, $false, $false, $false, $false, $true, 1, $false, "26+65=91", 2) | Out-Null
$d.Content.Find.Execute("97-91=6", $true, $false, $false, $false, $false, $true, 1, $false, "63-21=42", 2) | Out-Null
$d.Content.Find.Execute("45+41=86", $true, $false, $false, $false, $false, $true, 1, $false, "17+81=98", 2) | Out-Null
$d.Content.Find.Execute("67+2=69", $true, $false, $false, $false, $false, $true, 1, $false, "33-19=14", 2) | Out-Null
$d.Content.Find.Execute("64-12=52", $true, $false, $false, $false, $false, $true, 1, $false, "47+20=67", 2) | Out-Null
$d.Content.Find.Execute("92-78=14", $true, $false, $false, $false, $false, $true, 1, $false, "34+2=36", 2) | Out-Null
$d.Content.Find.Execute("48-34=14", $true, $false, $false, $false, $false, $true, 1, $false, "71+17=88", 2) | Out-Null
$d.Content.Find.Execute("36+51=87", $true, $false, $false, $false, $false, $true, 1, $false, "70-6=64", 2) | Out-Null
$d.Content.Find.Execute("67-6=61", $true, $false, $false, $false, $false, $true, 1, $false, "78-29=49", 2) | Out-Null
$d.Content.Find.Execute("61-35=26", $true, $false, $false, $false, $false, $true, 1, $false, "67+17=84", 2) | Out-Null
$d.Content.Find.Execute("42+19=61", $true, $false, $false, $false, $false, $true, 1, $false, "66-42=24", 2) | Out-Null
$d.Content.Find.Execute("80+2=82", $true, $false, $false, $false, $false, $true, 1, $false, "66+30=96", 2) | Out-Null
$d.Content.Find.Execute("99-35=64", $true, $false, $false, $false, $false, $true, 1, $false, "71-42=29", 2) | Out-Null
$d.Content.Find.Execute("40+11=51", $true, $false, $false, $false, $false, $true, 1, $false, "25+25=50", 2) | Out-Null
$d.Content.Find.Execute("52-22=30", $true, $false, $false, $false, $false, $true, 1, $false, "67-31=36", 2) | Out-Null
$d.Content.Find.Execute("30+16=46", $true, $false, $false, $false, $false, $true, 1, $false, "33-21=12", 2) | Out-Null
$d.Content.Find.Execute("82-80=2", $true, $false, $false, $false, $false, $true, 1, $false, "49-17=32", 2) | Out-Null
$d.Content.Find.Execute("25+43=68", $true, $false, $false, $false, $false, $true, 1, $false, "55+44=99", 2) | Out-Null
$d.Content.Find.Execute("46-41=5", $true, $false, $false, $false, $false, $true, 1, $false, "5+63=68", 2) | Out-Null
$d.Content.Find.Execute("50-40=10", $true, $false, $false, $false, $false, $true, 1, $false, "39+43=82", 2) | Out-Null
$d.Content.Find.Execute("33-4=29", $true, $false, $false, $false, $false, $true, 1, $false, "69-1=68", 2) | Out-Null
$d.Content.Find.Execute("44+7=51", $true, $false, $false, $false, $false, $true, 1, $false, "73-26=47", 2) | Out-Null
$d.Content.Find.Execute("47-9=38", $true, $false, $false, $false, $false, $true, 1, $false, "53-21=32", 2) | Out-Null
$d.Content.Find.Execute("73+14=87", $true, $false, $false, $false, $false, $true, 1, $false, "77-19=58", 2) | Out-Null
$d.Content.Find.Execute("44-32=12", $true, $false, $false, $false, $false, $true, 1, $false, "36-0=36", 2) | Out-Null
$d.Content.Find.Execute("79-29=50", $true, $false, $false, $false, $false, $true, 1, $false, "84-11=73", 2) | Out-Null
$d.Content.Find.Execute("17+66=83", $true, $false, $false, $false, $false, $true, 1, $false, "80-14=66", 2) | Out-Null
$d.Content.Find.Execute("75+1=76", $true, $false, $false, $false, $false, $true, 1, $false, "55-8=47", 2) | Out-Null
$d.Content.Find.Execute("70-5=65", $true, $false, $false, $false, $false, $true, 1, $false, "66+31=97", 2) | Out-Null
$d.Content.Find.Execute("19+1=20", $true, $false, $false, $false, $false, $true, 1, $false, "59-6=53", 2) | Out-Null
$d.Content.Find.Execute("20-8=12", $true, $false, $false, $false, $false, $true, 1, $false, "12+45=57", 2) | Out-Null
$d.Content.Find.Execute("29+21=50", $true, $false, $false, $false, $false, $true, 1, $false, "83-64=19", 2) | Out-Null
$d.Content.Find.Execute("96-85=11", $true, $false, $false, $false, $false, $true, 1, $false, "90-24=66", 2) | Out-Null
$d.Content.Find.Execute("88-1=87", $true, $false, $false, $false, $false, $true, 1, $false, "11+86=97", 2) | Out-Null
$d.Content.Find.Execute("52-3=49", $true, $false, $false, $false, $false, $true, 1, $false, "18+5=23", 2) | Out-Null
$d.Content.Find.Execute("43+51=94", $true, $false, $false, $false, $false, $true, 1, $false, "82-41=41", 2) | Out-Null
$d.Content.Find.Execute("95-15=80", $true, $false, $false, $false, $false, $true, 1, $false, "78-11=67", 2) | Out-Null
$d.Content.Find.Execute("99-40=59", $true, $false, $false, $false, $false, $true, 1, $false, "76-58=18", 2) | Out-Null
$d.Content.Find.Execute("85-43=42", $true, $false, $false, $false, $false, $true, 1, $false, "3+11=14", 2) | Out-Null
$d.Content.Find.Execute("29+19=48", $true, $false, $false, $false, $false, $true, 1, $false, "95-69=26", 2) | Out-Null
$d.Content.Find.Execute("17-12=5", $true, $false, $false, $false, $false, $true, 1, $false, "22-6=16", 2) | Out-Null
$d.Content.Find.Execute("61-43=18", $true, $false, $false, $false, $false, $true, 1, $false, "44+38=82", 2) | Out-Null
$d.Content.Find.Execute("86+7=93", $true, $false, $false, $false, $false, $true, 1, $false, "92-22=70", 2) | Out-Null
$d.Content.Find.Execute("93-75=18", $true, $false, $false, $false, $false, $true, 1, $false, "81-49=32", 2) | Out-Null
$d.Content.Find.Execute("13+48=61", $true, $false, $false, $false, $false, $true, 1, $false, "39+10=49", 2) | Out-Null
$d.Content.Find.Execute("8+2=10", $true, $false, $false, $false, $false, $true, 1, $false, "20+52=72", 2) | Out-Null
$d.Content.Find.Execute("82-31=51", $true, $false, $false, $false, $false, $true, 1, $false, "43+30=73", 2) | Out-Null
$d.Content.Find.Execute("34+51=85", $true, $false, $false, $false, $false, $true, 1, $false, "37-3=34", 2) | Out-Null
$d.Content.Find.Execute("74-61=13", $true, $false, $false, $false, $false, $true, 1, $false, "88-4=84", 2) | Out-Null
$d.Content.Find.Execute("1+47=48", $true, $false, $false, $false, $false, $true, 1, $false, "89+2=91", 2) | Out-Null
$d.Content.Find.Execute("37+2=39", $true, $false, $false, $false, $false, $true, 1, $false, "94-72=22", 2) | Out-Null
$d.Content.Find.Execute("28-22=6", $true, $false, $false, $false, $false, $true, 1, $false, "66-18=48", 2) | Out-Null
$d.Content.Find.Execute("34+32=66", $true, $false, $false, $false, $false, $true, 1, $false, "89+4=93", 2) | Out-Null
